{"js": "// Update the date paragraph (first paragraph of the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].load(\"text\");\nawait context.sync();\nif (paragraphs.items[0].text === \"2025-09-01 Monday\") {\n  paragraphs.items[0].insertText(\"2025-09-02 Tuesday\", Word.InsertLocation.replace);\n}\n\n// Update the division problems in the table (5 content rows, each followed\n// by 3 blank rows; 5 columns per row).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// cellUpdates[rowIndex] = [newCol0, newCol1, newCol2, newCol3, newCol4]\n// `null` means \"leave this cell untouched\".\nconst cellUpdates = {\n  0: [\"21\u00f77=\", \"19\u00f79=\", \"40\u00f75=\", \"76\u00f79=\", null],\n  4: [\"66\u00f73=\", \"58\u00f78=\", \"13\u00f74=\", \"38\u00f78=\", \"63\u00f74=\"],\n  8: [\"40\u00f72=\", \"71\u00f73=\", \"67\u00f79=\", \"53\u00f77=\", \"11\u00f72=\"],\n  12: [\"11\u00f74=\", \"98\u00f79=\", \"86\u00f74=\", \"18\u00f72=\", \"63\u00f78=\"],\n  16: [\"44\u00f74=\", \"71\u00f77=\", \"13\u00f76=\", \"46\u00f74=\", \"90\u00f75=\"],\n};\n\nfor (const rowIndexStr of Object.keys(cellUpdates)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = cellUpdates[rowIndex];\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < newValues.length; c++) {\n    if (newValues[c] !== null) {\n      cells.items[c].value = newValues[c];\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph of the document).\n$dateRange = $d.Paragraphs.Item(1).Range\nif ($dateRange.Text -like \"*2025-09-01 Monday*\") {\n    $dateRange.Text = \"2025-09-02 Tuesday\"\n}\n\n# Update the division problems in the table. The table has 20 rows: every\n# 4th row (1, 5, 9, 13, 17 in 1-based indexing) holds the 5 problem cells,\n# the rest are blank spacer rows.\n$table = $d.Tables.Item(1)\n\n$cellUpdates = @{\n    1  = @(\"21\u00f77=\", \"19\u00f79=\", \"40\u00f75=\", \"76\u00f79=\", $null)\n    5  = @(\"66\u00f73=\", \"58\u00f78=\", \"13\u00f74=\", \"38\u00f78=\", \"63\u00f74=\")\n    9  = @(\"40\u00f72=\", \"71\u00f73=\", \"67\u00f79=\", \"53\u00f77=\", \"11\u00f72=\")\n    13 = @(\"11\u00f74=\", \"98\u00f79=\", \"86\u00f74=\", \"18\u00f72=\", \"63\u00f78=\")\n    17 = @(\"44\u00f74=\", \"71\u00f77=\", \"13\u00f76=\", \"46\u00f74=\", \"90\u00f75=\")\n}\n\nforeach ($rowIndex in $cellUpdates.Keys) {\n    $newValues = $cellUpdates[$rowIndex]\n    for ($c = 0; $c -lt $newValues.Length; $c++) {\n        $newValue = $newValues[$c]\n        if ($null -ne $newValue) {\n            $table.Cell($rowIndex, $c + 1).Range.Text = $newValue\n        }\n    }\n}\n"}
